$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear out the old index column (A) -- its values are no longer used.
$ws.Range("A1:A7").Clear()

# 2. Insert a new blank row above the (former) header row; B:J shift down one row.
$ws.Range("A1").EntireRow.Insert()

Write-Host "structure done"
